$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)

$c = $g.GroupItems.Item("rc5")
$c.Left = (1234269/12700 + 0.00002)
$c.Top = (4708544/12700 + 0.00002)
$c.Width = (949086/12700 + 0.00002)
$c.Height = (822616/12700 + 0.00002)

$c = $g.GroupItems.Item("rc7")
$c.Left = (3343351/12700 + 0.00002)
$c.Top = (3827169/12700 + 0.00002)
$c.Width = (949086/12700 + 0.00002)
$c.Height = (1703992/12700 + 0.00002)

$c = $g.GroupItems.Item("rc8")
$c.Left = (4397892/12700 + 0.00002)
$c.Top = (4591028/12700 + 0.00002)
$c.Width = (949086/12700 + 0.00002)
$c.Height = (940133/12700 + 0.00002)

$c = $g.GroupItems.Item("rc9")
$c.Left = (5452433/12700 + 0.00002)
$c.Top = (4708544/12700 + 0.00002)
$c.Width = (949086/12700 + 0.00002)
$c.Height = (822616/12700 + 0.00002)

$c = $g.GroupItems.Item("rc10")
$c.Left = (6506974/12700 + 0.00002)
$c.Top = (3357102/12700 + 0.00002)
$c.Width = (949086/12700 + 0.00002)
$c.Height = (2174059/12700 + 0.00002)

$c = $g.GroupItems.Item("tx11")
$c.Left = (1608330/12700 + 0.00002)
$c.Top = (4123761/12700 + 0.00002)
$c.Width = (200965/12700 + 0.00002)
$c.Height = (129860/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "14"

$c = $g.GroupItems.Item("tx12")
$c.Left = (1518080/12700 + 0.00002)
$c.Top = (4344234/12700 + 0.00002)
$c.Width = (381464/12700 + 0.00002)
$c.Height = (169559/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "(9%)"

$c = $g.GroupItems.Item("tx13")
$c.Left = (2662871/12700 + 0.00002)
$c.Top = (2065013/12700 + 0.00002)
$c.Width = (200965/12700 + 0.00002)
$c.Height = (132065/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "49"

$c = $g.GroupItems.Item("tx14")
$c.TextFrame.TextRange.Text = "(31%)"

$c = $g.GroupItems.Item("tx15")
$c.Left = (3717412/12700 + 0.00002)
$c.Top = (3240180/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "29"

$c = $g.GroupItems.Item("tx16")
$c.Left = (3576921/12700 + 0.00002)
$c.Top = (3462858/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "(18%)"

$c = $g.GroupItems.Item("tx17")
$c.Left = (4771953/12700 + 0.00002)
$c.Top = (4004039/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "16"

$c = $g.GroupItems.Item("tx18")
$c.Left = (4631462/12700 + 0.00002)
$c.Top = (4226717/12700 + 0.00002)
$c.Width = (481947/12700 + 0.00002)
$c.Height = (169559/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "(10%)"

$c = $g.GroupItems.Item("tx19")
$c.Left = (5826494/12700 + 0.00002)
$c.Top = (4123761/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "14"

$c = $g.GroupItems.Item("tx20")
$c.Left = (5736245/12700 + 0.00002)
$c.Top = (4344234/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "(9%)"

$c = $g.GroupItems.Item("tx21")
$c.Left = (6881035/12700 + 0.00002)
$c.Top = (2770025/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "37"

$c = $g.GroupItems.Item("tx22")
$c.Left = (6740544/12700 + 0.00002)
$c.Top = (2992792/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "(23%)"

$c = $g.GroupItems.Item("tx25")
$c.Left = (824331/12700 + 0.00002)
$c.Top = (4299142/12700 + 0.00002)

$c = $g.GroupItems.Item("tx26")
$c.Left = (824331/12700 + 0.00002)
$c.Top = (3123975/12700 + 0.00002)

$c = $g.GroupItems.Item("tx27")
$c.Left = (824331/12700 + 0.00002)
$c.Top = (1948808/12700 + 0.00002)

$c = $g.GroupItems.Item("pl29")
$c.Left = (1031804/12700 + 0.00002)
$c.Top = (4355994/12700 + 0.00002)

$c = $g.GroupItems.Item("pl30")
$c.Left = (1031804/12700 + 0.00002)
$c.Top = (3180827/12700 + 0.00002)

$c = $g.GroupItems.Item("pl31")
$c.Left = (1031804/12700 + 0.00002)
$c.Top = (2005660/12700 + 0.00002)

$c = $g.GroupItems.Item("tx48")
$c.TextFrame.TextRange.Text = "159)"

$c = $g.GroupItems.Item("rc49")
$c.Left = (1234269/12700 + 0.00002)
$c.Top = (10201524/12700 + 0.00002)
$c.Width = (949086/12700 + 0.00002)
$c.Height = (746448/12700 + 0.00002)

$c = $g.GroupItems.Item("rc51")
$c.Left = (3343351/12700 + 0.00002)
$c.Top = (8868579/12700 + 0.00002)
$c.Width = (949086/12700 + 0.00002)
$c.Height = (2079392/12700 + 0.00002)

$c = $g.GroupItems.Item("rc52")
$c.Left = (4397892/12700 + 0.00002)
$c.Top = (10254841/12700 + 0.00002)
$c.Width = (949086/12700 + 0.00002)
$c.Height = (693130/12700 + 0.00002)

$c = $g.GroupItems.Item("rc53")
$c.Left = (5452433/12700 + 0.00002)
$c.Top = (10361477/12700 + 0.00002)
$c.Width = (949086/12700 + 0.00002)
$c.Height = (586495/12700 + 0.00002)

$c = $g.GroupItems.Item("rc54")
$c.Left = (6506974/12700 + 0.00002)
$c.Top = (9455075/12700 + 0.00002)
$c.Width = (949086/12700 + 0.00002)
$c.Height = (1492897/12700 + 0.00002)

$c = $g.GroupItems.Item("tx55")
$c.Left = (1608330/12700 + 0.00002)
$c.Top = (9616740/12700 + 0.00002)
$c.Width = (200965/12700 + 0.00002)
$c.Height = (129860/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "14"

$c = $g.GroupItems.Item("tx56")
$c.Left = (1518080/12700 + 0.00002)
$c.Top = (9837213/12700 + 0.00002)
$c.Width = (381464/12700 + 0.00002)
$c.Height = (169559/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "(9%)"

$c = $g.GroupItems.Item("tx58")
$c.TextFrame.TextRange.Text = "(34%)"

$c = $g.GroupItems.Item("tx59")
$c.Left = (3717412/12700 + 0.00002)
$c.Top = (8281502/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "39"

$c = $g.GroupItems.Item("tx60")
$c.Left = (3576921/12700 + 0.00002)
$c.Top = (8504269/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "(25%)"

$c = $g.GroupItems.Item("tx61")
$c.Left = (4771953/12700 + 0.00002)
$c.Top = (9667764/12700 + 0.00002)
$c.Width = (200965/12700 + 0.00002)
$c.Height = (132153/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "13"

$c = $g.GroupItems.Item("tx62")
$c.Left = (4681704/12700 + 0.00002)
$c.Top = (9890531/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "(8%)"

$c = $g.GroupItems.Item("tx63")
$c.Left = (5826494/12700 + 0.00002)
$c.Top = (9776693/12700 + 0.00002)
$c.Width = (200965/12700 + 0.00002)
$c.Height = (129860/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "11"

$c = $g.GroupItems.Item("tx64")
$c.Left = (5736245/12700 + 0.00002)
$c.Top = (9997166/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "(7%)"

$c = $g.GroupItems.Item("tx65")
$c.Left = (6881035/12700 + 0.00002)
$c.Top = (8868086/12700 + 0.00002)
$c.Width = (200965/12700 + 0.00002)
$c.Height = (132065/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "28"

$c = $g.GroupItems.Item("tx66")
$c.Left = (6740544/12700 + 0.00002)
$c.Top = (9090764/12700 + 0.00002)
$c.TextFrame.TextRange.Text = "(18%)"

$c = $g.GroupItems.Item("tx92")
$c.TextFrame.TextRange.Text = "159)"

